$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Relabel the battery/alarm current headers (F2/G2) for the "Add Panels" sheet.
$ws1.Range("F2").Value = "Alarm Current(A)"
$ws1.Range("G2").Value = "Standby Current(A)"

# Distinguish the RBus "Label" entry from the Gallery Type entry: K8 becomes
# "MPM800-1" (a new shared string) while I8 stays "MPM800".
$ws1.Range("K8").Value = "MPM800-1"

# Widen column K (Label) to fit the new text. The exported OOXML <col> width
# is ColumnWidth + 5/6, so 9.1666... (char units) round-trips to width="10".
$ws1.Columns.Item(11).ColumnWidth = 9.166666666666666

# Move the active selection to K8 (the cell that was just edited).
[void]$ws1.Range("K8").Select()
